{"js": "// Update the worksheet date and the 25 division problems to the new\n// \"three-digit number divided by one-digit number\" values.\nconst replacements = [\n  [\"2024-10-17 Thursday\", \"2024-10-18 Friday\"],\n  [\"585\u00f78=\", \"372\u00f72=\"],\n  [\"978\u00f79=\", \"827\u00f76=\"],\n  [\"655\u00f74=\", \"759\u00f78=\"],\n  [\"532\u00f72=\", \"751\u00f74=\"],\n  [\"415\u00f76=\", \"742\u00f79=\"],\n  [\"171\u00f79=\", \"762\u00f75=\"],\n  [\"339\u00f76=\", \"203\u00f73=\"],\n  [\"794\u00f73=\", \"967\u00f76=\"],\n  [\"540\u00f73=\", \"458\u00f72=\"],\n  [\"717\u00f75=\", \"420\u00f72=\"],\n  [\"133\u00f74=\", \"560\u00f72=\"],\n  [\"852\u00f76=\", \"399\u00f79=\"],\n  [\"702\u00f73=\", \"140\u00f78=\"],\n  [\"204\u00f78=\", \"500\u00f75=\"],\n  [\"955\u00f75=\", \"706\u00f77=\"],\n  [\"837\u00f76=\", \"240\u00f75=\"],\n  [\"576\u00f79=\", \"795\u00f77=\"],\n  [\"732\u00f72=\", \"429\u00f78=\"],\n  [\"826\u00f74=\", \"156\u00f79=\"],\n  [\"436\u00f76=\", \"633\u00f75=\"],\n  [\"988\u00f74=\", \"995\u00f79=\"],\n  [\"634\u00f77=\", \"899\u00f78=\"],\n  [\"603\u00f77=\", \"393\u00f74=\"],\n  [\"238\u00f73=\", \"137\u00f74=\"],\n  [\"160\u00f73=\", \"731\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 division problems to the new\n# \"three-digit number divided by one-digit number\" values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-10-17 Thursday\", \"2024-10-18 Friday\"),\n  @(\"585\u00f78=\", \"372\u00f72=\"),\n  @(\"978\u00f79=\", \"827\u00f76=\"),\n  @(\"655\u00f74=\", \"759\u00f78=\"),\n  @(\"532\u00f72=\", \"751\u00f74=\"),\n  @(\"415\u00f76=\", \"742\u00f79=\"),\n  @(\"171\u00f79=\", \"762\u00f75=\"),\n  @(\"339\u00f76=\", \"203\u00f73=\"),\n  @(\"794\u00f73=\", \"967\u00f76=\"),\n  @(\"540\u00f73=\", \"458\u00f72=\"),\n  @(\"717\u00f75=\", \"420\u00f72=\"),\n  @(\"133\u00f74=\", \"560\u00f72=\"),\n  @(\"852\u00f76=\", \"399\u00f79=\"),\n  @(\"702\u00f73=\", \"140\u00f78=\"),\n  @(\"204\u00f78=\", \"500\u00f75=\"),\n  @(\"955\u00f75=\", \"706\u00f77=\"),\n  @(\"837\u00f76=\", \"240\u00f75=\"),\n  @(\"576\u00f79=\", \"795\u00f77=\"),\n  @(\"732\u00f72=\", \"429\u00f78=\"),\n  @(\"826\u00f74=\", \"156\u00f79=\"),\n  @(\"436\u00f76=\", \"633\u00f75=\"),\n  @(\"988\u00f74=\", \"995\u00f79=\"),\n  @(\"634\u00f77=\", \"899\u00f78=\"),\n  @(\"603\u00f77=\", \"393\u00f74=\"),\n  @(\"238\u00f73=\", \"137\u00f74=\"),\n  @(\"160\u00f73=\", \"731\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
